# NIT-8040116241.xlsx: "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker row layout (columns B..J, rows 16..60) stays exactly the same
# (same worker, same style, same amounts). Only the "Periodo Mora" value in
# column E of each row is refreshed: the period list 1607..2003 used to run
# in ascending order down the rows (row16=1607 .. row60=2003); after the
# edit it is refreshed in descending order (row16=2003 .. row60=1607), i.e.
# the old periods are "removed" and the new ones are written in reverse.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
    "1712","1711","1710","1709","1708","1707","1706","1705","1704","1703","1702","1701",
    "1612","1611","1610","1609","1608","1607"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}
